$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")

# Clear out the sample times for the second device (column B, rows 2-8)
# so the downstream AVERAGE/MAX/MIN formulas on the Stats sheet that
# reference this column blow up with #DIV/0! instead of silently using
# stale/unassigned data (the delegate-style "unassigned -> NRE" bug being
# fixed here).
$data.Range("B2:B8").Clear()

# Switch the active sheet/selection to Data!B2 (previously Stats was the
# selected tab with Data!B9 selected).
$data.Activate()
$data.Range("B2").Select()
